# "Add files via upload" — the uploaded copy of the workbook drops the
# original title/header block (rows 3-4 of "Base de datos"), shifts the
# frequency table up so it starts at row 2, strips the ad-hoc cell
# styling/column widths/merges that went with the old header, turns the
# gridlines back on, and renames the sheet "Base de datos" -> "Hoja1"
# (sheetId 1 -> 2). Selection ends up on D6.

$wb  = $excel.ActiveWorkbook
$old = $wb.ActiveSheet

# Build the replacement sheet from scratch so it picks up Excel's plain
# defaults (no custom column widths, no merged cells, default page
# margins) instead of inheriting the old sheet's formatting.
$new = $wb.Worksheets.Add()
$new.Name = "Hoja1"

$new.Range("A2").Value = "Comunidad"
$new.Range("B2").Value = "Freq."

$new.Range("A3").Value = "Una vez a"
$new.Range("B3").Value = 35

$new.Range("A4").Value = "Una o dos"
$new.Range("B4").Value = 189

$new.Range("A5").Value = "Una o dos"
$new.Range("B5").Value = 323

$new.Range("A6").Value = "Nunca"
$new.Range("B6").Value = 929

$new.Range("A7").Value = "Total"
$new.Range("B7").Value = 1476

$new.Range("D6").Select()

# Drop the old, formatted sheet entirely.
$wb.Worksheets.Item("Base de datos").Delete()
